# Canada Weekly-Deaths-Prediction (r = 9): append new prediction rows.
#
# Rows 56-63 are an exact repeat of rows 48-55 (same week labels / values,
# same "KNN" model) - copy them down so the string/number cell types match
# exactly (avoids Excel's automatic text -> date coercion on a plain
# .Value assignment of "2021-01-09").
#
# Row 64 is a brand new week ("07 Mar -- 13 Mar 2021") with a new
# prediction value (62.15), same day-of-prediction (2021-01-09) and
# same model (KNN).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 56-63: duplicate rows 48-55 (columns A, B, D, F only - leave C/E
# untouched since the source rows don't populate them either).
$ws.Range("A48:B55").Copy($ws.Range("A56:B63"))
$ws.Range("D48:D55").Copy($ws.Range("D56:D63"))
$ws.Range("F48:F55").Copy($ws.Range("F56:F63"))

# Row 64: new week, new prediction value, same day/model as the rows above.
$ws.Range("A55:B55").Copy($ws.Range("A64:B64"))
$ws.Range("F55").Copy($ws.Range("F64"))
$ws.Range("B64").Value = "07 Mar -- 13 Mar 2021"
$ws.Range("D64").Value = 62.15
